$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, [string]$value) {
    $cell = $ws.Cells.Item($row, $col)
    if ($value -match '^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$') {
        $cell.NumberFormat = "@"
        $cell.Value = $value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $value
    }
}

# Row 2
Set-TextValue 2 4 "41.731.85"
Set-TextValue 2 5 "  +0.49%  "

# Row 3
Set-TextValue 3 4 "2.478.08"
Set-TextValue 3 5 "  +0.27%  "

# Row 4
Set-TextValue 4 5 "  +0.17%  "

# Row 5
Set-TextValue 5 4 "318.93"
Set-TextValue 5 5 "  +1.31%  "

# Row 6
Set-TextValue 6 4 "93.32"
Set-TextValue 6 5 "  +1.36%  "

# Row 7
Set-TextValue 7 5 "  +0.67%  "

# Row 8
Set-TextValue 8 5 "  +0.04%  "

# Row 9
Set-TextValue 9 4 "0.518"
Set-TextValue 9 5 "  +1.07%  "

# Row 10
Set-TextValue 10 2 "Avalanche"
Set-TextValue 10 3 "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue 10 4 "33.26"
Set-TextValue 10 5 "  +3.21%  "

# Row 11
Set-TextValue 11 2 "Dogecoin"
Set-TextValue 11 3 "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue 11 4 "0.0864"
Set-TextValue 11 5 "  +9.48%  "

# Row 12
Set-TextValue 12 5 "  +0.58%  "

# Row 13
Set-TextValue 13 4 "2.858.97"
Set-TextValue 13 5 "  +0.32%  "

# Row 14
Set-TextValue 14 4 "6.90"
Set-TextValue 14 5 "  +0.73%  "

# Row 15
Set-TextValue 15 4 "15.82"
Set-TextValue 15 5 "  -1.22%  "

# Row 16
Set-TextValue 16 4 "2.481.58"
Set-TextValue 16 5 "  +0.61%  "

# Row 17
Set-TextValue 17 4 "0.791"
Set-TextValue 17 5 "  +2.84%  "

# Row 18
Set-TextValue 18 4 "41.696.56"
Set-TextValue 18 5 "  +0.45%  "

# Row 19
Set-TextValue 19 5 "  -0.16%  "

# Row 20
Set-TextValue 20 4 "0.0₃0953"
Set-TextValue 20 5 "  +0.65%  "

# Row 21
Set-TextValue 21 4 "71.28"
Set-TextValue 21 5 "  -0.09%  "

# Row 22
Set-TextValue 22 4 "11.35"
Set-TextValue 22 5 "  +1.94%  "

# Row 23
Set-TextValue 23 4 "240.08"
Set-TextValue 23 5 "  +1.62%  "

# Row 24
Set-TextValue 24 5 "  +1.10%  "

# Row 25
Set-TextValue 25 4 "1.94"
Set-TextValue 25 5 "  +2.28%  "

# Row 26
Set-TextValue 26 5 "  +0.02%  "

# Row 27
Set-TextValue 27 4 "24.77"
Set-TextValue 27 5 "  +0.68%  "

# Row 28
Set-TextValue 28 4 "2.27"
Set-TextValue 28 5 "  +1.94%  "

# Row 29
Set-TextValue 29 4 "9.83"
Set-TextValue 29 5 "  +1.60%  "

# Row 30
Set-TextValue 30 4 "36.19"
Set-TextValue 30 5 "  +2.28%  "

# Row 31
Set-TextValue 31 4 "157.97"
Set-TextValue 31 5 "  +1.10%  "

# Row 32
Set-TextValue 32 4 "5.54"
Set-TextValue 32 5 "  +1.66%  "

# Row 33
Set-TextValue 33 5 "  -0.06%  "

# Row 34
Set-TextValue 34 2 "Hedera"
Set-TextValue 34 3 "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue 34 4 "0.0770"
Set-TextValue 34 5 "  +1.59%  "

# Row 35
Set-TextValue 35 2 "WEMIXToken"
Set-TextValue 35 3 "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue 35 4 "2.59"
Set-TextValue 35 5 "  +0.91%  "

# Row 36
Set-TextValue 36 4 "17.34"
Set-TextValue 36 5 "  +0.47%  "

# Row 37
Set-TextValue 37 4 "1.89"
Set-TextValue 37 5 "  +5.32%  "

# Row 38
Set-TextValue 38 4 "2.94"
Set-TextValue 38 5 "  +1.93%  "

# Row 39
Set-TextValue 39 5 "  +1.82%  "

# Row 40
Set-TextValue 40 5 "  +0.21%  "

# Row 41
Set-TextValue 41 4 "4.06"
Set-TextValue 41 5 "  +0.48%  "

# Row 42
Set-TextValue 42 5 "  +9.24%  "

# Row 43
Set-TextValue 43 2 "EnergySwap"
Set-TextValue 43 3 "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue 43 4 "19.54"
Set-TextValue 43 5 "  +5.63%  "

# Row 44
Set-TextValue 44 2 "Maker"
Set-TextValue 44 3 "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue 44 4 "1.997.18"
Set-TextValue 44 5 "  +2.76%  "

# Row 45
Set-TextValue 45 5 "  +1.02%  "

# Row 46
Set-TextValue 46 4 "3.00"
Set-TextValue 46 5 "  +2.34%  "

# Row 47
Set-TextValue 47 4 "9.34"
Set-TextValue 47 5 "  +3.13%  "

# Row 48
Set-TextValue 48 4 "2.716.20"
Set-TextValue 48 5 "  +0.38%  "

# Row 49
Set-TextValue 49 4 "97.54"
Set-TextValue 49 5 "  +0.43%  "

# Row 50
Set-TextValue 50 4 "74.47"
Set-TextValue 50 5 "  +3.72%  "

# Row 51
Set-TextValue 51 4 "67.37"
Set-TextValue 51 5 "  +0.19%  "

